$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.585.70"
$ws.Range("E2").Value = "  +2.58%  "
$ws.Range("D3").Value = "1.685.89"
$ws.Range("E3").Value = "  +3.38%  "
$ws.Range("E4").Value = "  -0.32%  "
$ws.Range("D5").Value = "'217.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.23%  "
$ws.Range("D6").Value = "'0.5352"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.21%  "
$ws.Range("E7").Value = "  -0.46%  "
$ws.Range("D8").Value = "'0.2681"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.42%  "
$ws.Range("D9").Value = "'0.06430"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.28%  "
$ws.Range("D10").Value = "'21.33"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.11%  "
$ws.Range("D11").Value = "'0.07760"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.00%  "
$ws.Range("D12").Value = "1.684.78"
$ws.Range("E12").Value = "  +3.39%  "
$ws.Range("D13").Value = "'4.498"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.83%  "
$ws.Range("D14").Value = "'0.5632"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.18%  "
$ws.Range("D15").Value = "0.0₅8396"
$ws.Range("E15").Value = "  +5.46%  "
$ws.Range("D16").Value = "'66.21"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.89%  "
$ws.Range("D17").Value = "26.586.48"
$ws.Range("E17").Value = "  +2.61%  "
$ws.Range("B18").Value = "Dai"
$ws.Range("C18").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D18").Value = "'1.000"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.41%  "
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").Value = "'4.825"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.02%  "
$ws.Range("D20").Value = "'194.94"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.53%  "
$ws.Range("D21").Value = "'10.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.66%  "
$ws.Range("D22").Value = "'6.402"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.95%  "
$ws.Range("D23").Value = "'1.002"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.32%  "
$ws.Range("D24").Value = "'143.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.91%  "
$ws.Range("D25").Value = "'0.1278"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.28%  "
$ws.Range("D26").Value = "'7.497"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.20%  "
$ws.Range("D27").Value = "'16.26"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.84%  "
$ws.Range("D28").Value = "'1.421"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.96%  "
$ws.Range("D29").Value = "'0.06133"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.36%  "
$ws.Range("D30").Value = "'1.279"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.08%  "
$ws.Range("D31").Value = "'3.607"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.56%  "
$ws.Range("D32").Value = "'3.464"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.71%  "
$ws.Range("D33").Value = "'1.706"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.31%  "
$ws.Range("E34").Value = "  +6.02%  "
$ws.Range("D35").Value = "'2.799"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.04%  "
$ws.Range("D36").Value = "'2.416"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.38%  "
$ws.Range("D37").Value = "'0.5731"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.90%  "
$ws.Range("D38").Value = "'0.01647"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.37%  "
$ws.Range("D39").Value = "'5.956"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.33%  "
$ws.Range("D40").Value = "'0.8692"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.53%  "
$ws.Range("D41").Value = "1.053.28"
$ws.Range("E41").Value = "  -1.68%  "
$ws.Range("E42").Value = "  -0.17%  "
$ws.Range("D43").Value = "'100.32"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.68%  "
$ws.Range("D44").Value = "1.835.84"
$ws.Range("E44").Value = "  +2.95%  "
$ws.Range("D45").Value = "'57.22"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.01%  "
$ws.Range("D46").Value = "0.0₈108"
$ws.Range("E46").Value = "  +1.41%  "
$ws.Range("D47").Value = "'8.188"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.42%  "
$ws.Range("D48").Value = "'1.001"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.20%  "
$ws.Range("D49").Value = "'6.103"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.64%  "
$ws.Range("D50").Value = "'0.05202"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.04%  "
$ws.Range("D51").Value = "'0.4239"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.35%  "
